$wb = $excel.ActiveWorkbook

# New "Ingles II" group-statistics row is inserted as row 11 (pushing the
# existing "4ARHV"/"6ARHV" summary rows down by one) on each of the three
# statistics sheets: "Estadisticos 1P", "Estadisticos 2P" and
# "Estadisticos Final". The "Rescatables" sheet needs no direct edit - its
# shared-string references stay meaningful automatically.

$sheetData = @{
    "Estadisticos 1P"    = @(11, 0, 1, 10, 90.91, 6.5)
    "Estadisticos 2P"    = @(11, 0, 0, 11, 100, 6.5)
    "Estadisticos Final" = @(11, 0, 0, 11, 100, 7.2)
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new blank row at position 11 - this shifts the old row 11
    # ("4ARHV" totals) down to row 12 and the old row 12 ("6ARHV" totals)
    # down to row 13, and grows the sheet's used range to A1:H13.
    $ws.Rows("11:11").Insert()

    $vals = $sheetData[$sheetName]

    $ws.Range("A11").Value = "Ingles II"
    $ws.Range("B11").Value = "2ASV"
    $ws.Range("C11").Value = $vals[0]
    $ws.Range("D11").Value = $vals[1]
    $ws.Range("E11").Value = $vals[2]
    $ws.Range("F11").Value = $vals[3]
    $ws.Range("G11").Value = $vals[4]
    $ws.Range("H11").Value = $vals[5]
}
